$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the columns that are going away (D: B_Additional_Benefits, E: B_Additional_Benefits_Format)
$ws.Range("D1:E3").Clear()

# Header row: A1 stays B_Id, B1 becomes B_Title, C1 becomes B_Limit (was column B)
$ws.Range("B1").Value = "B_Title"
$ws.Range("C1").Value = "B_Limit"

# Row 2: move numeric Limit from B2 into C2, put the name into B2
$ws.Range("C2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value = "mousa"

# Row 3: move numeric Limit from B3 into C3, put the name into B3
$ws.Range("C3").Value2 = $ws.Range("B3").Value2
$ws.Range("B3").Value = "sami"

# Column B width (new B column holding names) - matches column A's stored width of 5
$ws.Columns("B").ColumnWidth = 4.17

# sheet view / selection + sheet format changes
$ws.Range("G8").Select()
